$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").NumberFormat = "General"
$ws.Range("G2").NumberFormat = "General;General"
$ws.Range("G3").NumberFormat = "0_);[Red](0)"
Write-Host "done"
